# flower bottle support v3
# Duplicate Sheet1 into a new "v3" sheet, insert a new "center offset" row,
# tweak several inputs/formulas for the v3 revision, and append the new
# "cut rotation" row at the bottom.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# 1. Copy Sheet1 to a new sheet placed right after it, then rename it.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "v3"

# 2. Insert a new blank row at row 7 ("center offset"); this pushes every
#    row from the old row 7 downward by one and auto-adjusts formulas,
#    exactly like pressing Ctrl+"+" on a selected row in Excel.
$ws2.Rows.Item(7).Insert()

# 3. Updated input values for the v3 revision.
$ws2.Range("B2").Value = 35
$ws2.Range("B3").Value = 23
$ws2.Range("B6").Value = 2

# New row 7: "center offset" input (defaults to 0).
$ws2.Range("A7").Value = "center offset"
$ws2.Range("B7").Value = 0

# Row 8 ("slot min depth", shifted down from the old row 7) becomes 0.
$ws2.Range("B8").Value = 0

# "vial inner edge" now also folds in the new center-offset cell (B7).
$ws2.Range("C13").Formula = "=C2+B6+B8+B7"

# "vial outer edge reveal" flips sign.
$ws2.Range("C15").Value = -2

# "major diameter" also subtracts the new center-offset cell (B7).
$ws2.Range("C18").Formula = "=C14-C15-B7"

# 4. New standalone value below the existing block.
$ws2.Range("D20").Value = 89

# 5. New "cut rotation" row appended at the bottom.
$ws2.Range("A35").Value = "cut rotation"
$ws2.Range("B35").Value = 4
$ws2.Range("C35").Value = "degrees"

# 6. View/selection bookkeeping to mirror the authored state: Sheet1 keeps
#    a "select all, active cell E27" selection...
$ws1.Cells.Select()
$ws1.Range("E27").Activate()

# ...while v3 ends up the active, selected sheet with A36 selected.
$excel.Goto($ws2.Range("A36"), $true)
